$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Insert a new "Meta description" paragraph right after the title
#    (Heading 1) paragraph.
#
#    Target XML:
#      <w:p>
#        <w:r/>
#        <w:r><w:rPr><w:b/></w:rPr><w:t>Meta description</w:t></w:r>
#        <w:r><w:t>: Read our unbiased review of Blazin Hot 7s Stack Em Up
#          slot. Learn how to play the game and try it for free. Discover
#          pros and cons and RTP rate.</w:t></w:r>
#      </w:p>
#
#    A brand-new paragraph inserted right next to a styled (Heading)
#    paragraph inherits that heading style (and stamping the Style
#    property back to "Normal" leaves stray rsid attributes behind), so
#    instead we build the paragraph in a spot bordered by plain "Normal"
#    paragraphs (no explicit pStyle at all) and then relocate it with
#    cut/paste.
# ---------------------------------------------------------------------------

$metaBoldText = "Meta description"
$metaRestText = ": Read our unbiased review of Blazin Hot 7s Stack Em Up slot. Learn how to play the game and try it for free. Discover pros and cons and RTP rate."

$pGameplay = $d.Paragraphs.Item(3)
$pGameplay.Range.InsertParagraphBefore()
$pScratch = $d.Paragraphs.Item(3)

$scratchStart = $pScratch.Range.Start
$rText = $d.Range($scratchStart, $scratchStart)
$rText.Text = $metaBoldText + $metaRestText

$rBold = $d.Range($scratchStart, $scratchStart + $metaBoldText.Length)
$rBold.Bold = 1

# Re-fetch the (now filled-in) scratch paragraph and cut it as a whole
# (including its paragraph mark).
$pScratch = $d.Paragraphs.Item(3)
$cutRange = $d.Range($pScratch.Range.Start, $pScratch.Range.End)
$cutRange.Cut()

# Paste it back right after the title paragraph.
$pTitle = $d.Paragraphs.Item(1)
$pasteRange = $d.Range($pTitle.Range.End, $pTitle.Range.End)
$pasteRange.Paste()

$pMeta = $d.Paragraphs.Item(2)

# Prepend the leading empty run (<w:r/>) that is present in the target
# markup, using a raw OOXML insert so it lands *before* the text we just
# typed in rather than after it.
$metaStart = $pMeta.Range.Start
$rEmpty = $d.Range($metaStart, $metaStart)
$emptyRunXml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p>
            <w:r/>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
'@
$rEmpty.InsertXML($emptyRunXml)

# ---------------------------------------------------------------------------
# 2) Remove the duplicated bold-title paragraph near the end of the
#    document, and replace the italic paragraph's text that used to follow
#    it with the new feature-image prompt text.
# ---------------------------------------------------------------------------

$count = $d.Paragraphs.Count
$pDupTitle = $d.Paragraphs.Item($count - 1)
$delRange = $d.Range($pDupTitle.Range.Start, $pDupTitle.Range.End)
$delRange.Delete()

$count = $d.Paragraphs.Count
$pDescription = $d.Paragraphs.Item($count)

$descStart = $pDescription.Range.Start
$descEnd = $pDescription.Range.End - 1
$rDescription = $d.Range($descStart, $descEnd)

$newDescriptionText = 'Create a feature image fitting "Blazin Hot 7s Stack Em Up": - Draw a cartoon-style image of a happy Maya warrior with glasses wearing a headdress made of fruits such as cherries, oranges, lemons, plums, and watermelons. - Have the warrior holding a Stack''Em Up symbol in one hand and a handful of coins in the other hand. - Surround the warrior with cascading reels and colorful symbols. - Add text above the image that says "Blazin Hot 7s Stack Em Up" in bold, fiery letters.'
$rDescription.Text = $newDescriptionText

Write-Output "Done. Paragraph count: $($d.Paragraphs.Count)"
